$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the adjacency matrix with new GMM-derived edge weights.
# The matrix is symmetric, so each pair of cells (i,j) and (j,i) is updated together.

$ws.Range("T4").Value  = 0.153526745361367
$ws.Range("C21").Value = 0.153526745361367

$ws.Range("P5").Value  = 0.200166762977825
$ws.Range("D17").Value = 0.200166762977825

$ws.Range("M6").Value  = 0.123370294828301
$ws.Range("E14").Value = 0.123370294828301

$ws.Range("N7").Value  = 0.145100445147683
$ws.Range("F15").Value = 0.145100445147683
